# Generate Report for handoff
# Update "Latest Handoff Datetime" (column D) for the 7b164f77-... file row
# (row 5) on both the zh-cn and de-de sheets, recording a fresh handoff
# timestamp produced for the handback report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D5").Value = "2016-02-15 03:22:24"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D5").Value = "2016-02-15 03:22:37"
